$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.236.95'
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").Value = '2.027.28'
$ws.Range("E3").Value = '  -0.09%  '
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.08'
$ws.Range("E5").Value = '  +1.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.609'
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.32'
$ws.Range("E8").Value = '  +1.68%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.379'
$ws.Range("E9").Value = '  -1.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0785'
$ws.Range("E10").Value = '  -1.62%  '
$ws.Range("E11").Value = '  -2.28%  '
$ws.Range("D12").Value = '2.326.51'
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.34'
$ws.Range("E13").Value = '  -0.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.18'
$ws.Range("E14").Value = '  -2.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.741'
$ws.Range("E15").Value = '  -0.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.21'
$ws.Range("E16").Value = '  +0.51%  '
$ws.Range("D17").Value = '2.025.23'
$ws.Range("E17").Value = '  -0.17%  '
$ws.Range("D18").Value = '37.173.53'
$ws.Range("E18").Value = '  +0.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.15'
$ws.Range("E19").Value = '  -0.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.11'
$ws.Range("E20").Value = '  +0.49%  '
$ws.Range("D21").Value = '0.0₃0820'
$ws.Range("E21").Value = '  -2.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '223.28'
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("E24").Value = '  +1.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.23'
$ws.Range("E25").Value = '  -1.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.59'
$ws.Range("E26").Value = '  -2.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.07'
$ws.Range("E27").Value = '  -3.83%  '
$ws.Range("E28").Value = '  +2.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.76'
$ws.Range("E29").Value = '  -0.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.32'
$ws.Range("E30").Value = '  -1.26%  '
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.45'
$ws.Range("E32").Value = '  -0.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0602'
$ws.Range("E33").Value = '  -1.42%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.48'
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.96'
$ws.Range("E35").Value = '  +7.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.35'
$ws.Range("E36").Value = '  -0.80%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.21'
$ws.Range("E37").Value = '  +0.66%  '
$ws.Range("B38").Value = 'BinanceUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.51'
$ws.Range("E39").Value = '  +1.99%  '
$ws.Range("D40").Value = '1.473.20'
$ws.Range("E40").Value = '  -2.11%  '
$ws.Range("E41").Value = '  -2.29%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '94.42'
$ws.Range("E42").Value = '  -1.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.78'
$ws.Range("E43").Value = '  -2.26%  '
$ws.Range("E44").Value = '  -1.80%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.32'
$ws.Range("E45").Value = '  -2.92%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.21'
$ws.Range("E46").Value = '  +16.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.12'
$ws.Range("E47").Value = '  -1.80%  '
$ws.Range("E48").Value = '  +0.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.17'
$ws.Range("E49").Value = '  -0.36%  '
$ws.Range("E50").Value = '  +0.15%  '
$ws.Range("D51").Value = '2.213.72'
$ws.Range("E51").Value = '  +0.00%  '
